$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginData")

# Row 2 (TC_01 - "Check invalid login") test data password value updated
$ws.Range("C2").Value = "username=invalid@email.com,password=wer43345454"

# Row 2 expected result flips from Failure to Success
$ws.Range("D2").Value = "Success"

# Widen column C (TestData) independently of column B so it gets its own
# column width entry instead of being merged with column B's range
$ws.Columns.Item(3).ColumnWidth = 67.67
